{"js": "// Replace each two-digit multiplication problem's text with its new value.\n// The mapping below reflects the (old -> new) text changes from the diff,\n// applied in document order. Every old value is unique in the document, so\n// a simple search + replace per pair is safe (no cross-matches between\n// the old and new value sets either).\nconst replacements = [\n  [\"41\u00d793=\", \"63\u00d718=\"],\n  [\"79\u00d746=\", \"51\u00d773=\"],\n  [\"16\u00d795=\", \"62\u00d720=\"],\n  [\"73\u00d711=\", \"38\u00d788=\"],\n  [\"74\u00d785=\", \"57\u00d762=\"],\n  [\"21\u00d742=\", \"83\u00d750=\"],\n  [\"37\u00d747=\", \"30\u00d779=\"],\n  [\"55\u00d732=\", \"58\u00d720=\"],\n  [\"11\u00d727=\", \"69\u00d714=\"],\n  [\"17\u00d760=\", \"45\u00d723=\"],\n  [\"93\u00d786=\", \"21\u00d768=\"],\n  [\"71\u00d736=\", \"45\u00d789=\"],\n  [\"59\u00d712=\", \"97\u00d757=\"],\n  [\"68\u00d711=\", \"97\u00d764=\"],\n  [\"38\u00d745=\", \"29\u00d793=\"],\n  [\"19\u00d776=\", \"72\u00d719=\"],\n  [\"17\u00d732=\", \"13\u00d747=\"],\n  [\"27\u00d744=\", \"88\u00d797=\"],\n  [\"36\u00d735=\", \"59\u00d779=\"],\n  [\"33\u00d775=\", \"57\u00d755=\"],\n  [\"99\u00d750=\", \"95\u00d793=\"],\n  [\"17\u00d776=\", \"94\u00d774=\"],\n  [\"16\u00d759=\", \"33\u00d737=\"],\n  [\"37\u00d766=\", \"68\u00d777=\"],\n  [\"99\u00d787=\", \"17\u00d718=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem's text with its new value.\n# The pairs below reflect the (old -> new) text changes from the diff,\n# applied in document order. Every old value is unique in the document,\n# so a Find/Replace-all per pair is safe and unambiguous (it only ever\n# matches the single intended occurrence).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"41\u00d793=\", \"63\u00d718=\"),\n  @(\"79\u00d746=\", \"51\u00d773=\"),\n  @(\"16\u00d795=\", \"62\u00d720=\"),\n  @(\"73\u00d711=\", \"38\u00d788=\"),\n  @(\"74\u00d785=\", \"57\u00d762=\"),\n  @(\"21\u00d742=\", \"83\u00d750=\"),\n  @(\"37\u00d747=\", \"30\u00d779=\"),\n  @(\"55\u00d732=\", \"58\u00d720=\"),\n  @(\"11\u00d727=\", \"69\u00d714=\"),\n  @(\"17\u00d760=\", \"45\u00d723=\"),\n  @(\"93\u00d786=\", \"21\u00d768=\"),\n  @(\"71\u00d736=\", \"45\u00d789=\"),\n  @(\"59\u00d712=\", \"97\u00d757=\"),\n  @(\"68\u00d711=\", \"97\u00d764=\"),\n  @(\"38\u00d745=\", \"29\u00d793=\"),\n  @(\"19\u00d776=\", \"72\u00d719=\"),\n  @(\"17\u00d732=\", \"13\u00d747=\"),\n  @(\"27\u00d744=\", \"88\u00d797=\"),\n  @(\"36\u00d735=\", \"59\u00d779=\"),\n  @(\"33\u00d775=\", \"57\u00d755=\"),\n  @(\"99\u00d750=\", \"95\u00d793=\"),\n  @(\"17\u00d776=\", \"94\u00d774=\"),\n  @(\"16\u00d759=\", \"33\u00d737=\"),\n  @(\"37\u00d766=\", \"68\u00d777=\"),\n  @(\"99\u00d787=\", \"17\u00d718=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
